$wb = $excel.ActiveWorkbook

# --- Sheet "baseline" (sheet1): columns E and F get new values, G is a formula (E+F) that recalculates ---
$ws1 = $wb.Worksheets.Item("baseline")
$ws1.Cells.Item(2, 6).Value = 14
$ws1.Cells.Item(3, 6).Value = 14
$ws1.Cells.Item(4, 6).Value = 15
$ws1.Cells.Item(5, 6).Value = 16
$ws1.Cells.Item(6, 6).Value = 16
$ws1.Cells.Item(7, 6).Value = 17
$ws1.Cells.Item(8, 6).Value = 17
$ws1.Cells.Item(9, 6).Value = 17
$ws1.Cells.Item(10, 5).Value = 0
$ws1.Cells.Item(10, 6).Value = 18
$ws1.Cells.Item(11, 6).Value = 19
$ws1.Cells.Item(12, 6).Value = 20
$ws1.Cells.Item(13, 5).Value = 1
$ws1.Cells.Item(13, 6).Value = 21
$ws1.Cells.Item(14, 6).Value = 22
$ws1.Cells.Item(15, 6).Value = 23
$ws1.Cells.Item(16, 5).Value = 2
$ws1.Cells.Item(16, 6).Value = 24
$ws1.Cells.Item(17, 6).Value = 25
$ws1.Cells.Item(18, 5).Value = 3
$ws1.Cells.Item(18, 6).Value = 26
$ws1.Cells.Item(19, 5).Value = 4
$ws1.Cells.Item(19, 6).Value = 29
$ws1.Cells.Item(20, 5).Value = 4
$ws1.Cells.Item(20, 6).Value = 31
$ws1.Cells.Item(21, 5).Value = 5
$ws1.Cells.Item(21, 6).Value = 33

# --- Sheet "treejoin" (sheet2): columns F and G get new values, H is a formula (E+F+G) that recalculates ---
$ws2 = $wb.Worksheets.Item("treejoin")
$ws2.Cells.Item(2, 7).Value = 14
$ws2.Cells.Item(3, 7).Value = 14
$ws2.Cells.Item(4, 7).Value = 15
$ws2.Cells.Item(5, 7).Value = 16
$ws2.Cells.Item(6, 7).Value = 16
$ws2.Cells.Item(7, 7).Value = 17
$ws2.Cells.Item(8, 7).Value = 17
$ws2.Cells.Item(9, 6).Value = 0
$ws2.Cells.Item(9, 7).Value = 17
$ws2.Cells.Item(10, 6).Value = 0
$ws2.Cells.Item(10, 7).Value = 18
$ws2.Cells.Item(11, 6).Value = 0
$ws2.Cells.Item(11, 7).Value = 20
$ws2.Cells.Item(12, 6).Value = 0
$ws2.Cells.Item(12, 7).Value = 20
$ws2.Cells.Item(13, 6).Value = 0
$ws2.Cells.Item(13, 7).Value = 21
$ws2.Cells.Item(14, 6).Value = 0
$ws2.Cells.Item(14, 7).Value = 22
$ws2.Cells.Item(15, 6).Value = 0
$ws2.Cells.Item(15, 7).Value = 22
$ws2.Cells.Item(16, 6).Value = 0
$ws2.Cells.Item(16, 7).Value = 24
$ws2.Cells.Item(17, 6).Value = 0
$ws2.Cells.Item(17, 7).Value = 26
$ws2.Cells.Item(18, 6).Value = 0
$ws2.Cells.Item(18, 7).Value = 27
$ws2.Cells.Item(19, 6).Value = 0
$ws2.Cells.Item(19, 7).Value = 30
$ws2.Cells.Item(20, 6).Value = 0
$ws2.Cells.Item(20, 7).Value = 31
$ws2.Cells.Item(21, 6).Value = 1
$ws2.Cells.Item(21, 7).Value = 33

# --- Update the selections / top-left cell so the view matches the target workbook ---
$ws1.Activate()
$ws1.Range("A33").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2

$ws2.Activate()
$ws2.Range("G21").Select() | Out-Null
